# Apply data refresh to Top33_DataComp sheet (commit: "Tings are happening hard")
# Updates M2_Len/FX_Len counts and M2/FX first/last date serials for each country row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: China
$ws.Range("C2").Value = 342
$ws.Range("F2").Value = 45444
$ws.Range("G2").Value = 30319
$ws.Range("H2").Value = 45505

# Row 3: United States
$ws.Range("E3").Value = 30256
$ws.Range("F3").Value = 45444

# Row 4: Euro Area
$ws.Range("E4").Value = 30256
$ws.Range("F4").Value = 45444
$ws.Range("G4").Value = 30319
$ws.Range("H4").Value = 45505

# Row 5: Japan
$ws.Range("E5").Value = 30256
$ws.Range("F5").Value = 45444
$ws.Range("G5").Value = 30319
$ws.Range("H5").Value = 45505

# Row 6: United Kingdom
$ws.Range("C6").Value = 451
$ws.Range("F6").Value = 45444
$ws.Range("G6").Value = 30319
$ws.Range("H6").Value = 45505

# Row 7: South Korea
$ws.Range("E7").Value = 30256
$ws.Range("F7").Value = 45444
$ws.Range("G7").Value = 30319
$ws.Range("H7").Value = 45505

# Row 8: Hong Kong
$ws.Range("D8").Value = 418
$ws.Range("E8").Value = 30256
$ws.Range("F8").Value = 45444
$ws.Range("H8").Value = 45505

# Row 9: Australia
$ws.Range("E9").Value = 30256
$ws.Range("F9").Value = 45444
$ws.Range("G9").Value = 30319
$ws.Range("H9").Value = 45505

# Row 10: Taiwan
$ws.Range("D10").Value = 491
$ws.Range("E10").Value = 30256
$ws.Range("F10").Value = 45444
$ws.Range("H10").Value = 45505

# Row 11: Canada
$ws.Range("E11").Value = 30225
$ws.Range("F11").Value = 45413
$ws.Range("G11").Value = 30319
$ws.Range("H11").Value = 45505

# Row 12: Russia
$ws.Range("C12").Value = 379
$ws.Range("D12").Value = 360
$ws.Range("F12").Value = 45444
$ws.Range("H12").Value = 45505

# Row 13: Switzerland
$ws.Range("C13").Value = 474
$ws.Range("F13").Value = 45413
$ws.Range("G13").Value = 30319
$ws.Range("H13").Value = 45505

# Row 14: Brazil
$ws.Range("D14").Value = 404
$ws.Range("H14").Value = 45505

# Row 15: India
$ws.Range("G15").Value = 30286
$ws.Range("H15").Value = 45505

# Row 16: Mexico
$ws.Range("C16").Value = 462
$ws.Range("D16").Value = 418
$ws.Range("F16").Value = 45413
$ws.Range("H16").Value = 45505

# Row 17: Saudi Arabia
$ws.Range("C17").Value = 378
$ws.Range("D17").Value = 402
$ws.Range("F17").Value = 45444
$ws.Range("H17").Value = 45505

# Row 18: Singapore
$ws.Range("D18").Value = 268
$ws.Range("E18").Value = 30256
$ws.Range("F18").Value = 45444
$ws.Range("H18").Value = 45505

# Row 19: Indonesia
$ws.Range("D19").Value = 406
$ws.Range("E19").Value = 30256
$ws.Range("F19").Value = 45444
$ws.Range("H19").Value = 45505

# Row 20: Malaysia
$ws.Range("C20").Value = 489
$ws.Range("F20").Value = 45444
$ws.Range("G20").Value = 30319
$ws.Range("H20").Value = 45505

# Row 21: Sweden
$ws.Range("C21").Value = 317
$ws.Range("F21").Value = 45413
$ws.Range("G21").Value = 30319
$ws.Range("H21").Value = 45505

# Row 22: Poland
$ws.Range("C22").Value = 331
$ws.Range("D22").Value = 375
$ws.Range("F22").Value = 45444
$ws.Range("H22").Value = 45505

# Row 23: Israel
$ws.Range("D23").Value = 323
$ws.Range("E23").Value = 30225
$ws.Range("F23").Value = 45413
$ws.Range("H23").Value = 45505

# Row 24: Egypt
$ws.Range("C24").Value = 341
$ws.Range("D24").Value = 321
$ws.Range("F24").Value = 45444
$ws.Range("H24").Value = 45505

# Row 25: Norway
$ws.Range("E25").Value = 30225
$ws.Range("F25").Value = 45413
$ws.Range("G25").Value = 30319
$ws.Range("H25").Value = 45505

# Row 26: Philippines
$ws.Range("D26").Value = 388
$ws.Range("E26").Value = 30256
$ws.Range("F26").Value = 45444
$ws.Range("H26").Value = 45505

# Row 27: New Zealand
$ws.Range("D27").Value = 225
$ws.Range("E27").Value = 30225
$ws.Range("F27").Value = 45413
$ws.Range("H27").Value = 45505

# Row 28: Denmark
$ws.Range("C28").Value = 402
$ws.Range("F28").Value = 45444
$ws.Range("G28").Value = 30319
$ws.Range("H28").Value = 45505

# Row 29: South Africa
$ws.Range("G29").Value = 30319
$ws.Range("H29").Value = 45505

# Row 30: Chile
$ws.Range("C30").Value = 462
$ws.Range("D30").Value = 406
$ws.Range("F30").Value = 45444
$ws.Range("H30").Value = 45505

# Row 31: Bangladesh
$ws.Range("C31").Value = 412
$ws.Range("D31").Value = 325
$ws.Range("F31").Value = 45413
$ws.Range("H31").Value = 45505

# Row 32: Colombia
$ws.Range("D32").Value = 418
$ws.Range("E32").Value = 30256
$ws.Range("F32").Value = 45444
$ws.Range("H32").Value = 45505

# Row 33: Morocco
$ws.Range("C33").Value = 474
$ws.Range("D33").Value = 325
$ws.Range("F33").Value = 45444
$ws.Range("H33").Value = 45505

# Row 34: Kuwait
$ws.Range("C34").Value = 367
$ws.Range("D34").Value = 375
$ws.Range("F34").Value = 45444
$ws.Range("H34").Value = 45505
